$d = $word.ActiveDocument

$replacements = @(
    @{old = "87×21=1827"; new = "94×46=4324"},
    @{old = "56×72=4032"; new = "78×30=2340"},
    @{old = "14×22=308"; new = "34×41=1394"},
    @{old = "40×75=3000"; new = "64×89=5696"},
    @{old = "66×15=990"; new = "12×88=1056"},
    @{old = "56×35=1960"; new = "73×55=4015"},
    @{old = "24×67=1608"; new = "29×61=1769"},
    @{old = "59×13=767"; new = "41×87=3567"},
    @{old = "91×34=3094"; new = "65×60=3900"},
    @{old = "72×51=3672"; new = "16×93=1488"},
    @{old = "64×54=3456"; new = "76×17=1292"},
    @{old = "34×65=2210"; new = "74×92=6808"},
    @{old = "82×32=2624"; new = "68×99=6732"},
    @{old = "13×79=1027"; new = "89×51=4539"},
    @{old = "31×86=2666"; new = "94×35=3290"},
    @{old = "37×14=518"; new = "48×11=528"},
    @{old = "23×93=2139"; new = "99×62=6138"},
    @{old = "44×75=3300"; new = "18×23=414"},
    @{old = "33×47=1551"; new = "80×98=7840"},
    @{old = "52×50=2600"; new = "66×48=3168"},
    @{old = "20×78=1560"; new = "40×81=3240"},
    @{old = "23×87=2001"; new = "30×89=2670"},
    @{old = "77×55=4235"; new = "96×55=5280"},
    @{old = "97×22=2134"; new = "98×21=2058"},
    @{old = "54×67=3618"; new = "80×14=1120"}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
